$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3072
$ws.Range("I64").Value = 3108.3333
$ws.Range("J64").Value = 2999.3333
$ws.Range("K64").Value = 3108.3333
$ws.Range("L64").Value = 2999.3333
$ws.Range("M64").Value = -2860.3333
$ws.Range("N64").Value = -3495.3333

$ws.Range("H67").Value = 3072
$ws.Range("I67").Value = 3108.3333
$ws.Range("J67").Value = 2999.3333
$ws.Range("K67").Value = 3108.3333
$ws.Range("L67").Value = 2999.3333
$ws.Range("M67").Value = -2250.3333
$ws.Range("N67").Value = -4715.3333

$ws.Range("H69").Value = 1512.6
$ws.Range("I69").Value = 1512.6
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 4537.799999999999
$ws.Range("L69").Value = 0
$ws.Range("M69").ClearContents()
$ws.Range("N69").Value = -3663.799999999999

$ws.Range("H72").Value = 1512.6
$ws.Range("I72").Value = 1512.6
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 13613.4
$ws.Range("L72").Value = 0
$ws.Range("M72").ClearContents()
$ws.Range("N72").Value = -9245.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 763.5238000000001
$ws.Range("I74").Value = 596.7059
$ws.Range("J74").Value = 1472.5
$ws.Range("K74").Value = 596.7059
$ws.Range("L74").Value = 1472.5
$ws.Range("M74").Value = 277.2941
$ws.Range("N74").Value = -3220.5

$ws.Range("H77").Value = 763.5238000000001
$ws.Range("I77").Value = 596.7059
$ws.Range("J77").Value = 1472.5
$ws.Range("K77").Value = 2983.5295
$ws.Range("L77").Value = 7362.5
$ws.Range("M77").Value = 1384.4705
$ws.Range("N77").Value = -16098.5

$ws.Range("H132").Value = 3049.0881
$ws.Range("I132").Value = 2234.7083
$ws.Range("J132").Value = 5003.6
$ws.Range("K132").Value = 6704.124899999999
$ws.Range("L132").Value = 15010.8
$ws.Range("M132").Value = -4174.124899999999
$ws.Range("N132").Value = -20070.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 70014.5
$ws.Range("I21").Value = 70013
$ws.Range("J21").Value = 70015
$ws.Range("K21").Value = 70013
$ws.Range("L21").Value = 70015
$ws.Range("M21").Value = -69778
$ws.Range("N21").Value = -70485

$ws.Range("H62").Value = 5298.3335
$ws.Range("I62").Value = 4205
$ws.Range("J62").Value = 5845
$ws.Range("K62").Value = 4205
$ws.Range("L62").Value = 5845
$ws.Range("M62").Value = -3581
$ws.Range("N62").Value = -7093

$ws.Range("H65").Value = 5298.3335
$ws.Range("I65").Value = 4205
$ws.Range("J65").Value = 5845
$ws.Range("K65").Value = 21025
$ws.Range("L65").Value = 29225
$ws.Range("M65").Value = -17905
$ws.Range("N65").Value = -35465

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H55").Value = 1417.1111
$ws.Range("I55").Value = 600.8
$ws.Range("J55").Value = 2437.5
$ws.Range("K55").Value = 1802.4
$ws.Range("L55").Value = 7312.5
$ws.Range("M55").Value = -1625.4
$ws.Range("N55").Value = -7666.5

$ws.Range("H64").Value = 12041.7
$ws.Range("I64").Value = 2000.5
$ws.Range("J64").Value = 14552
$ws.Range("K64").Value = 6001.5
$ws.Range("L64").Value = 43656
$ws.Range("M64").Value = -5731.5
$ws.Range("N64").Value = -44196

$ws.Range("H67").Value = 12041.7
$ws.Range("I67").Value = 2000.5
$ws.Range("J67").Value = 14552
$ws.Range("K67").Value = 6001.5
$ws.Range("L67").Value = 43656
$ws.Range("M67").Value = -5065.5
$ws.Range("N67").Value = -45528

$ws.Range("H70").Value = 3333.1667
$ws.Range("I70").Value = 2000
$ws.Range("J70").Value = 3999.75
$ws.Range("K70").Value = 6000
$ws.Range("L70").Value = 11999.25
$ws.Range("M70").Value = -5685
$ws.Range("N70").Value = -12629.25

$ws.Range("H73").Value = 3333.1667
$ws.Range("I73").Value = 2000
$ws.Range("J73").Value = 3999.75
$ws.Range("K73").Value = 6000
$ws.Range("L73").Value = 11999.25
$ws.Range("M73").Value = -4908
$ws.Range("N73").Value = -14183.25

$ws.Range("H114").Value = 906.9524
$ws.Range("I114").Value = 219.33333
$ws.Range("J114").Value = 1021.55554
$ws.Range("K114").Value = 657.99999
$ws.Range("L114").Value = 3064.66662
$ws.Range("M114").Value = 2596.00001
$ws.Range("N114").Value = -9572.66662

$ws.Range("H117").Value = 1310.125
$ws.Range("I117").Value = 195.5
$ws.Range("J117").Value = 1469.3572
$ws.Range("K117").Value = 586.5
$ws.Range("L117").Value = 4408.071599999999
$ws.Range("M117").Value = 2855.5
$ws.Range("N117").Value = -11292.0716

$ws.Range("H121").Value = 26324.25
$ws.Range("I121").Value = 287.5
$ws.Range("J121").Value = 39342.625
$ws.Range("K121").Value = 862.5
$ws.Range("L121").Value = 118027.875
$ws.Range("M121").Value = 447.5
$ws.Range("N121").Value = -120647.875

$ws.Range("H130").Value = 2500
$ws.Range("J130").Value = 2500
$ws.Range("L130").Value = 7500
$ws.Range("N130").Value = -17540

$ws.Range("H131").Value = 2146.2593
$ws.Range("J131").Value = 1802.1666
$ws.Range("L131").Value = 5406.4998
$ws.Range("N131").Value = -15486.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 28650
$ws.Range("J68").Value = 28650
$ws.Range("L68").Value = 28650
$ws.Range("N68").Value = -30272

$ws.Range("H71").Value = 28650
$ws.Range("J71").Value = 28650
$ws.Range("L71").Value = 85950
$ws.Range("N71").Value = -94062

$ws.Range("H110").Value = 40000
$ws.Range("J110").Value = 40000
$ws.Range("L110").Value = 40000
$ws.Range("N110").Value = -48180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1330.3462
$ws.Range("I46").Value = 549
$ws.Range("J46").Value = 1472.409
$ws.Range("K46").Value = 549
$ws.Range("L46").Value = 1472.409
$ws.Range("M46").Value = -361
$ws.Range("N46").Value = -1848.409

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 8465.843999999999
$ws.Range("I132").Value = 1993.1702
$ws.Range("J132").Value = 18606.367
$ws.Range("K132").Value = 5979.5106
$ws.Range("L132").Value = 55819.101
$ws.Range("M132").Value = -3449.5106
$ws.Range("N132").Value = -60879.101
